{"js": "// The document contains 6 occurrences of an <id> tag whose numeric value\n// is currently split across three differently-formatted runs, e.g.\n//   \"<id>\" + \"p083r_a1\" + \"</id>\"\n// The edit collapses each of these into a single run (using the\n// formatting of the first/third run) containing the fully-resolved text,\n// and renumbers the id by dropping the \"a\" prefix, e.g.\n//   \"<id>p083r_1</id>\"\n// Do this for ids 1 through 6.\nconst body = context.document.body;\n\nfor (let i = 1; i <= 6; i++) {\n  const oldText = `<id>p083r_a${i}</id>`;\n  const newText = `<id>p083r_${i}</id>`;\n\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    // Replacing the whole matched range with the final text merges the\n    // three runs into a single run, inheriting the formatting of the\n    // range's first run (the Courier-New styled \"<id>\" run).\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# The document contains 6 occurrences of an <id> tag whose numeric value\n# is currently split across three differently-formatted runs, e.g.\n#   \"<id>\" + \"p083r_a1\" + \"</id>\"\n# The edit collapses each of these into a single run (using the\n# formatting of the first/third run) containing the fully-resolved text,\n# and renumbers the id by dropping the \"a\" prefix, e.g.\n#   \"<id>p083r_1</id>\"\n# Do this for ids 1 through 6.\n$d = $word.ActiveDocument\n\nfor ($i = 1; $i -le 6; $i++) {\n    $oldText = \"<id>p083r_a$i</id>\"\n    $newText = \"<id>p083r_$i</id>\"\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    # Replacing the whole matched range merges the three runs into a\n    # single run, inheriting the formatting of the range's first run\n    # (the Courier-New styled \"<id>\" run).\n    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
